$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (matching the original inline-string
# storage) instead of being auto-converted to numbers by Excel while we set them.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '26.791.01'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '1.640.07'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').Value = '218.44'
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').Value = '0.502'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('E7').Value = '  -0.51%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('D9').Value = '0.0622'
$ws.Range('E9').Value = '  -0.89%  '
$ws.Range('D10').Value = '19.21'
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('D11').Value = '0.0846'
$ws.Range('D12').Value = '1.869.50'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').Value = '1.639.75'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').Value = '4.15'
$ws.Range('E14').Value = '  -0.73%  '
$ws.Range('D15').Value = '0.526'
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').Value = '64.98'
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('D17').Value = '26.808.23'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('E18').Value = '  -0.92%  '
$ws.Range('D19').Value = '216.29'
$ws.Range('E19').Value = '  +0.95%  '
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').Value = '6.66'
$ws.Range('E21').Value = '  +5.90%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '4.36'
$ws.Range('E22').Value = '  -0.26%  '
$ws.Range('D23').Value = '2.34'
$ws.Range('E23').Value = '  -2.62%  '
$ws.Range('D24').Value = '9.16'
$ws.Range('E24').Value = '  -2.32%  '
$ws.Range('D25').Value = '147.47'
$ws.Range('E25').Value = '  +1.41%  '
$ws.Range('E26').Value = '  -0.44%  '
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('D28').Value = '7.10'
$ws.Range('E28').Value = '  +0.11%  '
$ws.Range('D29').Value = '15.72'
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').Value = '0.0507'
$ws.Range('E30').Value = '  -1.47%  '
$ws.Range('D31').Value = '1.20'
$ws.Range('E31').Value = '  +1.13%  '
$ws.Range('E32').Value = '  +1.78%  '
$ws.Range('E33').Value = '  -0.96%  '
$ws.Range('E34').Value = '  +0.82%  '
$ws.Range('D35').Value = '1.264.63'
$ws.Range('E35').Value = '  -2.19%  '
$ws.Range('D36').Value = '2.45'
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').Value = '0.0175'
$ws.Range('E37').Value = '  -0.56%  '
$ws.Range('D38').Value = '0.530'
$ws.Range('E38').Value = '  -1.57%  '
$ws.Range('D39').Value = '0.817'
$ws.Range('E39').Value = '  -1.31%  '
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('D41').Value = '0.806'
$ws.Range('E41').Value = '  -0.41%  '
$ws.Range('E42').Value = '  -0.31%  '
$ws.Range('D43').Value = '1.780.20'
$ws.Range('E43').Value = '  -0.81%  '
$ws.Range('E44').Value = '  -4.34%  '
$ws.Range('D45').Value = '92.31'
$ws.Range('E45').Value = '  +0.90%  '
$ws.Range('D46').Value = '60.67'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').Value = '1.59'
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.0516'
$ws.Range('E48').Value = '  -0.91%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '0.0963'
$ws.Range('E49').Value = '  -1.76%  '
$ws.Range('E50').Value = '  -2.09%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '1.01'
$ws.Range('E51').Value = '  -0.44%  '

# Restore the default (Normal) style on column D so no visible formatting changes
# beyond the cell values themselves.
$priceRange.Style = "Normal"

